$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "touch" stat entry to "work" (UI_WorkItem feature)
$ws.Range("B3").Value = "work"

# Update the active cell selection to B4
$ws.Range("B4").Select()
